{"js": "// Remove the stray `_GoBack` bookmark that originally sat at the very start\n// of the document (around the \"Ejercicio Gen\u00e9rico\" heading). Deleting it\n// first frees up bookmark id \"0\" so the new bookmark inserted below reuses\n// that same id, matching the target OOXML exactly.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst body = context.document.body;\n\n// --- \"\u00bfEs P\" -> \"\u00bfEs \" + italic \"P\" -------------------------------------\nconst pHits = body.search(\"\u00bfEs P\", { matchCase: true });\nawait context.sync();\n\nif (pHits.items.length > 0) {\n  const pHit = pHits.items[0];\n  const pLetter = pHit.search(\"P\", { matchCase: true });\n  await context.sync();\n  if (pLetter.items.length > 0) {\n    pLetter.items[0].font.italic = true;\n  }\n}\n\n// --- \"\u00bfEs M = {Mam\u00edferos...\" -> \"\u00bfEs \" + italic \"M\" (wrapped in the\n//     relocated _GoBack bookmark) + \" = {Mam\u00edferos...}\" ------------------\nconst mHits = body.search(\"\u00bfEs M\", { matchCase: true });\nawait context.sync();\n\nif (mHits.items.length > 0) {\n  const mHit = mHits.items[0];\n  const mLetter = mHit.search(\"M\", { matchCase: true });\n  await context.sync();\n  if (mLetter.items.length > 0) {\n    const mRange = mLetter.items[0];\n    mRange.font.italic = true;\n    mRange.insertBookmark(\"_GoBack\");\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Remove the stray `_GoBack` bookmark that originally sat at the very start\n# of the document (around the \"Ejercicio Gen\u00e9rico\" heading). Deleting it\n# first frees up bookmark id \"0\" so the new bookmark inserted below reuses\n# that same id, matching the target OOXML exactly.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# --- \"\u00bfEs P\" -> \"\u00bfEs \" + italic \"P\" -------------------------------------\n$find1 = $d.Content\n$find1.Find.Text = \"\u00bfEs P\"\n$find1.Find.MatchCase = $true\n$found1 = $find1.Find.Execute()\nif ($found1) {\n    $pLetter = $d.Range($find1.End - 1, $find1.End)\n    $pLetter.Font.Italic = $true\n}\n\n# --- \"\u00bfEs M = {Mam\u00edferos...\" -> \"\u00bfEs \" + italic \"M\" (wrapped in the\n#     relocated _GoBack bookmark) + \" = {Mam\u00edferos...}\" ------------------\n$find2 = $d.Content\n$find2.Find.Text = \"\u00bfEs M\"\n$find2.Find.MatchCase = $true\n$found2 = $find2.Find.Execute()\nif ($found2) {\n    $mLetter = $d.Range($find2.End - 1, $find2.End)\n    $mLetter.Font.Italic = $true\n    $d.Bookmarks.Add(\"_GoBack\", $mLetter)\n}\n"}
